$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: values sourced from original row 12
$ws.Range("D2").Value = 44466
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 2000

# Row 3: values sourced from original row 2
$ws.Range("D3").Value = 44473
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 180
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 2000

# Row 4: values sourced from original row 11
$ws.Range("D4").Value = 44434
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 2000

# Row 5: values sourced from original row 6
$ws.Range("D5").Value = 44503
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("S5").Value = 3000

# Row 6: values sourced from original row 7
$ws.Range("D6").Value = 44503
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("S6").Value = 2500

# Row 7: values sourced from original row 8
$ws.Range("D7").Value = 44511
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 28000
$ws.Range("O7").Value = 28000
$ws.Range("P7").Value = 28000
$ws.Range("S7").Value = 2800

# Row 8: values sourced from original row 10
$ws.Range("D8").Value = 44432
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 2000

# Row 9: unchanged (maps to itself)
# Row 10: values sourced from original row 5
$ws.Range("D10").Value = 44435
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 2000

# Row 11: values sourced from original row 3
$ws.Range("D11").Value = 44517
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 27000
$ws.Range("O11").Value = 27000
$ws.Range("P11").Value = 27000
$ws.Range("S11").Value = 2700

# Row 12: values sourced from original row 4
$ws.Range("D12").Value = 44517
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 25000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 25000
$ws.Range("S12").Value = 2500

